$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.808.95'
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").Value = '1.909.44'
$ws.Range("E3").Value = '  +0.92%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '312.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5167'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +5.99%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3779'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.43%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07249'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.94%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '21.33'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.67%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.9053'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.955.50'
$ws.Range("E12").Value = '  +3.44%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07656'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.17%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.455'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.31%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '92.13'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.96%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.10%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.000008715'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.46%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").Value = '27.845.02'
$ws.Range("E19").Value = '  +0.41%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.53'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.54%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.152'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").Value = '2.202.17'
$ws.Range("E22").Value = '  +3.01%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.85'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.636'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.26%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '153.74'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.35'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.167'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.49%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '114.98'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.12%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.857'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09092'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.04%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.861'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +4.98%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.181'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.50%  '

$ws.Range("E34").Value = '  +0.80%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.7789'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.70%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.02094'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.77%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.603'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.86%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.077'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.14%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.5593'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("E40").Value = '  -0.16%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.716'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.56%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '115.39'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.11%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '8.572'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.81%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.1517'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4828'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.93%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.51'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.90%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.9982'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.13%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.618'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.13%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '66.96'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.72%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05991'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.12%  '
